# Cleaned Test data to have valid and invalid credentials in their respective sheets
#
# Sheet1 ends up holding only VALID credentials (admin@yourstore.com / admin,
# repeated for every row) and Sheet2 ends up holding INVALID credentials
# (wrong username/password combinations), including a brand new
# user3@test.com / pass126 row.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")

# ---------------------------------------------------------------------------
# Sheet1: fix up rows 4-5 to the valid admin/admin pair, then drop the old
# rows 6-7 (which are no longer needed once every row is a valid login) and
# the stale mailto hyperlinks that pointed at the removed rows.
# ---------------------------------------------------------------------------
$ws1.Range("B4").Value = "admin"
$ws1.Range("A5").Value = "admin@yourstore.com"
$ws1.Rows("6:7").Delete()
$ws1.Hyperlinks.Delete()

# ---------------------------------------------------------------------------
# Sheet2: replace row 5 with a new (invalid) user3@test.com / pass126 login
# and wire up its mailto hyperlink, preserving the existing cell style.
# ---------------------------------------------------------------------------
$sheet2A5Style = $ws2.Range("A5").Style
$ws2.Range("A5").Value = "user3@test.com"
$ws2.Range("B5").Value = "pass126"
$ws2.Hyperlinks.Add($ws2.Range("A5"), "mailto:user3@test.com")
$ws2.Range("A5").Style = $sheet2A5Style

# ---------------------------------------------------------------------------
# Selection / active sheet bookkeeping to match the saved view state:
# Sheet2's cursor rests on B5, and Sheet1 (with A5:B5 selected) is the
# active tab shown when the workbook re-opens.
# ---------------------------------------------------------------------------
$ws2.Activate() | Out-Null
$ws2.Range("B5").Select() | Out-Null

$ws1.Activate() | Out-Null
$ws1.Range("A5:B5").Select() | Out-Null
